# Updated cryptos list - applies the latest price/volume snapshot to the
# cryptocurrency tracking worksheet (Sheet1). Mirrors the scheduled
# GitHub Actions refresh: updates the "Price" (column D) and
# "Volume(1h)" (column E) values for each coin row.
#
# Some "Price" values are plain decimal-looking numbers (e.g. "214.85").
# The source sheet stores ALL of column D/E as literal text (so values such
# as "1.00" keep their trailing zero and values like "25.875.94", which use
# a dotted-thousands style, stay intact instead of being treated as dates).
# To stop Excel's automatic "looks like a number" conversion from mangling
# those cells, we briefly mark them as Text before assigning, then restore
# the default "Normal" cell style so no formatting residue is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

$ws.Range('D2').Value = '25.875.94'
$ws.Range('E2').Value = '  +0.54%  '
$ws.Range('D3').Value = '1.635.43'
$ws.Range('E3').Value = '  +0.06%  '
$ws.Range('E4').Value = '  -0.25%  '
Set-TextValue 'D5' '214.85'
$ws.Range('E5').Value = '  -0.26%  '
Set-TextValue 'D6' '0.503'
$ws.Range('E6').Value = '  -0.34%  '
Set-TextValue 'D7' '1.00'
$ws.Range('E7').Value = '  -0.29%  '
$ws.Range('E8').Value = '  -0.77%  '
$ws.Range('E9').Value = '  -0.67%  '
Set-TextValue 'D10' '19.64'
$ws.Range('E10').Value = '  +0.13%  '
$ws.Range('E11').Value = '  +0.57%  '
Set-TextValue 'D12' '4.26'
$ws.Range('E12').Value = '  +0.47%  '
$ws.Range('D13').Value = '1.861.18'
$ws.Range('E13').Value = '  +0.06%  '
$ws.Range('D14').Value = '1.637.65'
$ws.Range('E14').Value = '  +0.08%  '
Set-TextValue 'D15' '0.552'
$ws.Range('E15').Value = '  -0.88%  '
$ws.Range('E16').Value = '  -0.77%  '
Set-TextValue 'D17' '62.89'
$ws.Range('E17').Value = '  +0.21%  '
$ws.Range('D18').Value = '25.863.26'
$ws.Range('E18').Value = '  +0.36%  '
$ws.Range('E19').Value = '  -0.23%  '
$ws.Range('E20').Value = '  -0.08%  '
Set-TextValue 'D21' '191.64'
$ws.Range('E21').Value = '  -1.09%  '
Set-TextValue 'D22' '9.98'
$ws.Range('E22').Value = '  +0.29%  '
Set-TextValue 'D23' '6.33'
$ws.Range('E23').Value = '  +0.63%  '
Set-TextValue 'D24' '0.999'
$ws.Range('E24').Value = '  -0.32%  '
Set-TextValue 'D25' '1.80'
Set-TextValue 'D26' '142.44'
$ws.Range('E26').Value = '  +0.51%  '
$ws.Range('E27').Value = '  +1.00%  '
$ws.Range('E28').Value = '  -0.52%  '
$ws.Range('E29').Value = '  -0.11%  '
$ws.Range('E30').Value = '  -0.49%  '
Set-TextValue 'D31' '0.0494'
$ws.Range('E31').Value = '  +0.57%  '
$ws.Range('E32').Value = '  +0.15%  '
$ws.Range('E33').Value = '  -0.33%  '
$ws.Range('E34').Value = '  +0.69%  '
$ws.Range('E35').Value = '  +0.39%  '
Set-TextValue 'D36' '0.907'
$ws.Range('E36').Value = '  +0.69%  '
$ws.Range('D37').Value = '1.148.99'
$ws.Range('E37').Value = '  +2.24%  '
$ws.Range('E38').Value = '  -0.69%  '
$ws.Range('E39').Value = '  -0.93%  '
$ws.Range('E40').Value = '  +0.70%  '
$ws.Range('E41').Value = '  -0.26%  '
Set-TextValue 'D42' '5.63'
$ws.Range('E42').Value = '  +0.97%  '
Set-TextValue 'D43' '100.58'
$ws.Range('E43').Value = '  +0.83%  '
Set-TextValue 'D44' '0.801'
$ws.Range('E44').Value = '  -0.13%  '
$ws.Range('D45').Value = '1.771.14'
$ws.Range('E45').Value = '  +0.09%  '
$ws.Range('D46').Value = '0.0₆0110'
$ws.Range('E46').Value = '  -1.54%  '
Set-TextValue 'D47' '55.62'
$ws.Range('E47').Value = '  +0.90%  '
$ws.Range('E49').Value = '  +5.37%  '
$ws.Range('E50').Value = '  +0.01%  '
Set-TextValue 'D51' '7.57'
$ws.Range('E51').Value = '  +0.42%  '
